$d = $word.ActiveDocument

$pairs = @(
    @("548×3=", "189×9="),
    @("439×4=", "267×4="),
    @("987×7=", "761×6="),
    @("483×4=", "973×6="),
    @("486×2=", "283×5="),
    @("878×8=", "991×9="),
    @("983×2=", "958×6="),
    @("281×5=", "698×4="),
    @("671×3=", "901×2="),
    @("282×2=", "713×2="),
    @("852×4=", "911×8="),
    @("952×2=", "729×6="),
    @("597×4=", "131×2="),
    @("977×5=", "231×6="),
    @("977×3=", "769×3="),
    @("280×5=", "463×4="),
    @("709×6=", "598×5="),
    @("549×8=", "969×8="),
    @("932×9=", "444×9="),
    @("685×5=", "910×5="),
    @("320×9=", "492×7="),
    @("907×7=", "429×9="),
    @("637×2=", "396×3="),
    @("477×3=", "555×4="),
    @("995×3=", "295×5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
